# The "2" sheet (first tab, tabSelected) holds the query-table data that
# feeds the table `_2` (A1:B25). The edit re-orders the bin1..bin17 /
# bins / frequency / id / max / mean / min / size rows into a natural
# (numeric-aware) ascending order by Name, bumps bin17's value from 3 to 4,
# and adds a running-total formula in C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-write column A (Name) / column B (Value) in the new row order.
$ws.Range("A2").Value = "bin1"
$ws.Range("B2").Value = 74

$ws.Range("A3").Value = "bin2"
$ws.Range("B3").Value = 53

$ws.Range("A4").Value = "bin3"
$ws.Range("B4").Value = 48

$ws.Range("A5").Value = "bin4"
$ws.Range("B5").Value = 35

$ws.Range("A6").Value = "bin5"
$ws.Range("B6").Value = 25

$ws.Range("A7").Value = "bin6"
$ws.Range("B7").Value = 18

$ws.Range("A8").Value = "bin7"
$ws.Range("B8").Value = 10

$ws.Range("A9").Value = "bin8"
$ws.Range("B9").Value = 12

$ws.Range("A10").Value = "bin9"
$ws.Range("B10").Value = 7

$ws.Range("A11").Value = "bin10"
$ws.Range("B11").Value = 3

$ws.Range("A12").Value = "bin11"
$ws.Range("B12").Value = 1

$ws.Range("A13").Value = "bin12"
$ws.Range("B13").Value = 5

$ws.Range("A14").Value = "bin13"
$ws.Range("B14").Value = 2

$ws.Range("A15").Value = "bin14"
$ws.Range("B15").Value = 3

$ws.Range("A16").Value = "bin15"
$ws.Range("B16").Value = 0

$ws.Range("A17").Value = "bin16"
$ws.Range("B17").Value = 0

$ws.Range("A18").Value = "bin17"
$ws.Range("B18").Value = 4

$ws.Range("A19").Value = "bins"
$ws.Range("B19").Value = 17

$ws.Range("A20").Value = "frequency"
$ws.Range("B20").Value = 6.1169411764705881

$ws.Range("A21").Value = "id"
$ws.Range("B21").Value = 2

$ws.Range("A22").Value = "max"
$ws.Range("B22").Value = 104.01900000000001

$ws.Range("A23").Value = "mean"
$ws.Range("B23").Value = 20.632756666666669

$ws.Range("A24").Value = "min"
$ws.Range("B24").Value = 0.031

$ws.Range("A25").Value = "size"
$ws.Range("B25").Value = 300

# New running-total column: sum of the bin1..bin17 frequencies.
$ws.Range("C1").Formula = "=SUM(B2:B18)"

# Matches the author's final selection/active cell in the saved file.
$ws.Range("C1").Select() | Out-Null
